$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.962.73"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.559.60"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.00%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'208.08"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  +0.12%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'22.13"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D12").Value = "1.781.12"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "1.558.48"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E14").Value = "  +0.12%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.521"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +0.24%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'61.92"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "26.956.76"
$ws.Range("D18").Value = "0.0₃0709"
$ws.Range("E18").Value = "  +1.81%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'216.33"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  +0.40%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'1.94"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.65%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'152.83"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +0.03%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.0474"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").Value = "1.427.86"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("E36").Value = "  +8.66%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'2.33"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("E38").Value = "  +1.27%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.534"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +2.69%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'5.87"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +2.97%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.809"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("E44").Value = "  +0.47%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'64.66"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "1.694.28"
$ws.Range("E47").Value = "  +0.07%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'87.37"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  +5.74%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  +0.10%  "
